$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (Octopus-Inspired Soft Hydrogel Robots) - add funding text (ONR award)
$ws.Range("E9").Value = "This  work  is  supported  by  Office  of  Naval  Research  Award  N00014-17-1-2117"

# Row 19 - new project entry: SCRAM Platform 1
$ws.Range("A19").Value = "SCRAM Platform 1"
$ws.Range("B19").Value = "Mohammad Sharifzadeh, Yuhao Jiang"
$ws.Range("C19").Value = "coming soon…"
$ws.Range("E19").Value = 'This work is supported by <a href="https://www.nsf.gov/awardsearch/showAward?AWD_ID=1935324">NSF Award #1935324</a>'

# Row 8 (Digging Robot) - add funding text (NSF Award #1841574)
$ws.Range("E8").Value = 'This work is supported by <a href="https://www.nsf.gov/awardsearch/showAward?AWD_ID=1841574">NSF Award #1841574</a>'

# Update the active selection to match the new cursor position
$ws.Range("E8").Select()
